$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.642.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.597.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.582.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.624.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0512"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.274.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("E35").Value = "  -11.08%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +18.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.836"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.734.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.101"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "

Write-Host "Updated cryptos list values"
